$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 384
$ws.Range("I38").Value = 384
$ws.Range("K38").Value = 1152
$ws.Range("M38").Value = -780
$ws.Range("H41").Value = 1729.1852
$ws.Range("I41").Value = 1634.7
$ws.Range("J41").Value = 1999.1428
$ws.Range("K41").Value = 1634.7
$ws.Range("L41").Value = 1999.1428
$ws.Range("M41").Value = -1194.7
$ws.Range("N41").Value = -2879.1428
$ws.Range("H62").Value = 76497.36
$ws.Range("I62").Value = 95332.45
$ws.Range("K62").Value = 95332.45
$ws.Range("M62").Value = -94708.45
$ws.Range("H65").Value = 76497.36
$ws.Range("I65").Value = 95332.45
$ws.Range("K65").Value = 476662.25
$ws.Range("M65").Value = -473542.25
$ws.Range("H100").Value = 2596.8235
$ws.Range("I100").Value = 2366.818
$ws.Range("K100").Value = 2366.818
$ws.Range("M100").Value = -1825.818
$ws.Range("H107").Value = 683.8823
$ws.Range("I107").Value = 558.9286
$ws.Range("J107").Value = 1267
$ws.Range("K107").Value = 558.9286
$ws.Range("L107").Value = 1267
$ws.Range("M107").Value = 1361.0714
$ws.Range("N107").Value = -5107
$ws.Range("H135").Value = 39364.77
$ws.Range("I135").Value = 821.3333
$ws.Range("K135").Value = 7391.9997
$ws.Range("M135").Value = -4856.9997
$ws.Range("H138").Value = 2106.3076
$ws.Range("I138").Value = 1453.4286
$ws.Range("J138").Value = 2548.5806
$ws.Range("K138").Value = 4360.2858
$ws.Range("L138").Value = 7645.7418
$ws.Range("M138").Value = 779.7142000000003
$ws.Range("N138").Value = -17925.7418
$ws.Range("H141").Value = 31870
$ws.Range("I141").Value = 31870
$ws.Range("K141").Value = 95610
$ws.Range("M141").Value = -90430

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 52698.61
$ws.Range("I32").Value = 29017.459
$ws.Range("K32").Value = 29017.459
$ws.Range("M32").Value = -28730.459
$ws.Range("H74").Value = 2296.2964
$ws.Range("I74").Value = 2042.875
$ws.Range("K74").Value = 2042.875
$ws.Range("M74").Value = -1168.875
$ws.Range("H77").Value = 2296.2964
$ws.Range("I77").Value = 2042.875
$ws.Range("K77").Value = 10214.375
$ws.Range("M77").Value = -5846.375
$ws.Range("H110").Value = 3901.6155
$ws.Range("I110").Value = 3810.0833
$ws.Range("K110").Value = 3810.0833
$ws.Range("M110").Value = -1765.0833
$ws.Range("H122").Value = 2322.7778
$ws.Range("I122").Value = 2297.4666
$ws.Range("J122").Value = 2449.3333
$ws.Range("K122").Value = 6892.399800000001
$ws.Range("L122").Value = 7347.999899999999
$ws.Range("M122").Value = -4442.399800000001
$ws.Range("N122").Value = -12247.9999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5760.946
$ws.Range("J20").Value = 8616.875
$ws.Range("L20").Value = 8616.875
$ws.Range("N20").Value = -9110.875
$ws.Range("H86").Value = 2311.6667
$ws.Range("I86").Value = 2007.5555
$ws.Range("J86").Value = 2767.8333
$ws.Range("K86").Value = 2007.5555
$ws.Range("L86").Value = 2767.8333
$ws.Range("M86").Value = -884.5554999999999
$ws.Range("N86").Value = -5013.8333
$ws.Range("H89").Value = 2311.6667
$ws.Range("I89").Value = 2007.5555
$ws.Range("J89").Value = 2767.8333
$ws.Range("K89").Value = 10037.7775
$ws.Range("L89").Value = 13839.1665
$ws.Range("M89").Value = -4421.7775
$ws.Range("N89").Value = -25071.1665
$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1255.8889
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1329
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1329
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1903
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = ""
$ws.Range("H113").Value = 1255.8889
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1329
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1329
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5669

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 533
$ws.Range("I5").Value = 424.5
$ws.Range("J5").Value = 750
$ws.Range("K5").Value = 1273.5
$ws.Range("L5").Value = 2250
$ws.Range("M5").Value = -1161.5
$ws.Range("N5").Value = -2474
$ws.Range("H12").Value = 500.85715
$ws.Range("J12").Value = 550.3
$ws.Range("L12").Value = 1650.9
$ws.Range("N12").Value = -1996.9
$ws.Range("H34").Value = 3135.1
$ws.Range("I34").Value = 190
$ws.Range("J34").Value = 6080.2
$ws.Range("K34").Value = 570
$ws.Range("L34").Value = 18240.6
$ws.Range("M34").Value = -486
$ws.Range("N34").Value = -18408.6
$ws.Range("H97").Value = 2647.7778
$ws.Range("I97").Value = 1750
$ws.Range("J97").Value = 2760
$ws.Range("K97").Value = 5250
$ws.Range("L97").Value = 8280
$ws.Range("M97").Value = -4754
$ws.Range("N97").Value = -9272
$ws.Range("H98").Value = 824
$ws.Range("J98").Value = 496.5
$ws.Range("L98").Value = 1489.5
$ws.Range("N98").Value = -4485.5
$ws.Range("H107").Value = 829.0571
$ws.Range("I107").Value = 574.25
$ws.Range("J107").Value = 904.55554
$ws.Range("K107").Value = 1722.75
$ws.Range("L107").Value = 2713.66662
$ws.Range("M107").Value = 197.25
$ws.Range("N107").Value = -6553.66662
$ws.Range("H122").Value = 940.1905
$ws.Range("J122").Value = 983.8
$ws.Range("L122").Value = 8854.199999999999
$ws.Range("N122").Value = -13754.2
$ws.Range("H135").Value = 533
$ws.Range("I135").Value = 424.5
$ws.Range("J135").Value = 750
$ws.Range("K135").Value = 3820.5
$ws.Range("L135").Value = 6750
$ws.Range("M135").Value = -1285.5
$ws.Range("N135").Value = -11820

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5852.1113
$ws.Range("I70").Value = 5663.2666
$ws.Range("J70").Value = 6796.3335
$ws.Range("K70").Value = 5663.2666
$ws.Range("L70").Value = 6796.3335
$ws.Range("M70").Value = -5393.2666
$ws.Range("N70").Value = -7336.3335
$ws.Range("H73").Value = 5852.1113
$ws.Range("I73").Value = 5663.2666
$ws.Range("J73").Value = 6796.3335
$ws.Range("K73").Value = 5663.2666
$ws.Range("L73").Value = 6796.3335
$ws.Range("M73").Value = -4727.2666
$ws.Range("N73").Value = -8668.333500000001
$ws.Range("H132").Value = 2438.0527
$ws.Range("I132").Value = 2519
$ws.Range("J132").Value = 1750
$ws.Range("K132").Value = 7557
$ws.Range("L132").Value = 5250
$ws.Range("M132").Value = -5027
$ws.Range("N132").Value = -10310

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1099.75
$ws.Range("I22").Value = 966.3333
$ws.Range("K22").Value = 966.3333
$ws.Range("M22").Value = -671.3333
$ws.Range("H27").Value = 1099.75
$ws.Range("I27").Value = 966.3333
$ws.Range("K27").Value = 966.3333
$ws.Range("M27").Value = -859.3333
$ws.Range("H92").Value = 30387.2
$ws.Range("I92").Value = 30380
$ws.Range("K92").Value = 30380
$ws.Range("M92").Value = -27884
$ws.Range("H129").Value = 90495
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 90495
$ws.Range("K129").Value = 0
$ws.Range("M129").Value = ""
$ws.Range("N129").Value = -100495
$ws.Range("L129").Value = 90495
$ws.Range("H132").Value = 2450.8975
$ws.Range("J132").Value = 7348.8
$ws.Range("L132").Value = 22046.4
$ws.Range("N132").Value = -27106.4

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 23248.666
$ws.Range("J63").Value = 23248.666
$ws.Range("L63").Value = 23248.666
$ws.Range("N63").Value = -24496.666
$ws.Range("H66").Value = 23248.666
$ws.Range("J66").Value = 23248.666
$ws.Range("L66").Value = 69745.99800000001
$ws.Range("N66").Value = -75985.99800000001
$ws.Range("H122").Value = 2310.2856
$ws.Range("I122").Value = 1682.9286
$ws.Range("K122").Value = 5048.7858
$ws.Range("M122").Value = -2598.7858
$ws.Range("H132").Value = 7418.5557
$ws.Range("I132").Value = 8308.875
$ws.Range("J132").Value = 7418.5557
$ws.Range("K132").Value = 24926.625
$ws.Range("M132").Value = -22396.625
